$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "287.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-7.35%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.52%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.034"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.09%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07311"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-4.92%"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.13%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.547"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-8.56%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9086"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.91%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1199"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.98%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1754"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.41%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08646"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.13%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04158"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1051"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.04%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001278"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.09%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005746"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.39%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.399"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.32%"

$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.397"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.16%"

$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3284"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.09%"

$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.563"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.18%"

$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1342"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.73%"

$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2886"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.05%"

$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.03840"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.68%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001268"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003668"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-13.63%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001284"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.96%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003728"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02326"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-8.30%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.04990"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.01%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007704"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.69%"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "162.98%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.30%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007378"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.99%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007538"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.66%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3105"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.50%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006520"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.84%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.03%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "15.54%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.03%"
